# Populate the header row + 5 data rows with the new dropdown-option data
# (Order_ID / Transaksi / Kategori / Keterangan / Status / Order_By table).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Order_ID"
$ws.Range("B1").Value = "Transaksi"
$ws.Range("C1").Value = "Kategori"
$ws.Range("D1").Value = "Keterangan"
$ws.Range("E1").Value = "Status"
$ws.Range("F1").Value = "Order_By"

$ws.Range("A2").Value = "WO0123594"
$ws.Range("B2").Value = "MO"
$ws.Range("C2").Value = "INDIHOME"
$ws.Range("D2").Value = "Order ke 1"
$ws.Range("E2").Value = "Order"
$ws.Range("F2").Value = "Plasa"

$ws.Range("A3").Value = "WO0123590"
$ws.Range("B3").Value = "DO"
$ws.Range("C3").Value = "INDIBIZ"
$ws.Range("D3").Value = "Order ke 2"
$ws.Range("E3").Value = "Pickup"
$ws.Range("F3").Value = "Teknisi"

$ws.Range("A4").Value = "WO0123591"
$ws.Range("B4").Value = "RO"
$ws.Range("C4").Value = "WIFIID"
$ws.Range("D4").Value = "Order ke 3"
$ws.Range("E4").Value = "Close"
$ws.Range("F4").Value = "Plasa"

$ws.Range("A5").Value = "WO0123592"
$ws.Range("B5").Value = "PSB"
$ws.Range("C5").Value = "OLO"
$ws.Range("D5").Value = "Order ke 4"
$ws.Range("E5").Value = "Order"
$ws.Range("F5").Value = "Teknisi"

$ws.Range("A6").Value = "WO0123593"
$ws.Range("B6").Value = "PDA"
$ws.Range("C6").Value = "VPNIP"
$ws.Range("D6").Value = "Order ke 5"
$ws.Range("E6").Value = "CLose"
$ws.Range("F6").Value = "Plasa"

# Column widths (characters) tuned so the saved OOXML width lands on the
# same pixel width the source workbook uses for each column.
$ws.Columns("A").ColumnWidth = 16.5
$ws.Columns("B").ColumnWidth = 15.333333333333334
$ws.Columns("C").ColumnWidth = 20.166666666666668
$ws.Columns("D").ColumnWidth = 19.166666666666668
$ws.Columns("E").ColumnWidth = 16.833333333333332
$ws.Columns("F").ColumnWidth = 17.666666666666668

# Thin border around every cell of the table.
$tbl = $ws.Range("A1:F6")
$tbl.Borders.LineStyle = 1
$tbl.Borders.Weight = 2

# Status column (E) gets a white (Background 1 theme) solid fill on top of
# the border.
$ws.Range("E1:E6").Interior.ThemeColor = 2

# Leave the active selection where the author left it.
$ws.Range("I8").Select()
